$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.504.23'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('D2').Style = 'Normal'

$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +5.41%  '
$ws.Range('E2').NumberFormat = 'General'
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.725.58'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +4.69%  '
$ws.Range('E3').NumberFormat = 'General'
$ws.Range('E3').Style = 'Normal'

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('E4').NumberFormat = 'General'
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.59'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'

$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.33%  '
$ws.Range('E5').NumberFormat = 'General'
$ws.Range('E5').Style = 'Normal'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5349'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.90%  '
$ws.Range('E6').NumberFormat = 'General'
$ws.Range('E6').Style = 'Normal'

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('E7').NumberFormat = 'General'
$ws.Range('E7').Style = 'Normal'

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2666'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('D8').Style = 'Normal'

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.49%  '
$ws.Range('E8').NumberFormat = 'General'
$ws.Range('E8').Style = 'Normal'

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06579'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('D9').Style = 'Normal'

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +4.37%  '
$ws.Range('E9').NumberFormat = 'General'
$ws.Range('E9').Style = 'Normal'

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.58'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('D10').Style = 'Normal'

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +6.21%  '
$ws.Range('E10').NumberFormat = 'General'
$ws.Range('E10').Style = 'Normal'

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07702'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.25%  '
$ws.Range('E11').NumberFormat = 'General'
$ws.Range('E11').Style = 'Normal'

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.605'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('D12').Style = 'Normal'

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.35%  '
$ws.Range('E12').NumberFormat = 'General'
$ws.Range('E12').Style = 'Normal'

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.726.87'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('D13').Style = 'Normal'

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +4.74%  '
$ws.Range('E13').NumberFormat = 'General'
$ws.Range('E13').Style = 'Normal'

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.963.85'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = 'Normal'

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +4.74%  '
$ws.Range('E14').NumberFormat = 'General'
$ws.Range('E14').Style = 'Normal'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5817'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('D15').Style = 'Normal'

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +4.44%  '
$ws.Range('E15').NumberFormat = 'General'
$ws.Range('E15').Style = 'Normal'

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8273'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D16').Style = 'Normal'

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.96%  '
$ws.Range('E16').NumberFormat = 'General'
$ws.Range('E16').Style = 'Normal'

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.80'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('D17').Style = 'Normal'

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +4.16%  '
$ws.Range('E17').NumberFormat = 'General'
$ws.Range('E17').Style = 'Normal'

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '27.516.58'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('D18').Style = 'Normal'

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +5.52%  '
$ws.Range('E18').NumberFormat = 'General'
$ws.Range('E18').Style = 'Normal'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '217.79'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('D19').Style = 'Normal'

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +13.04%  '
$ws.Range('E19').NumberFormat = 'General'
$ws.Range('E19').Style = 'Normal'

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('E20').NumberFormat = 'General'
$ws.Range('E20').Style = 'Normal'

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.725'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.73%  '
$ws.Range('E21').NumberFormat = 'General'
$ws.Range('E21').Style = 'Normal'

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.50%  '
$ws.Range('E22').NumberFormat = 'General'
$ws.Range('E22').Style = 'Normal'

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.81%  '
$ws.Range('E23').NumberFormat = 'General'
$ws.Range('E23').Style = 'Normal'

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.17%  '
$ws.Range('E24').NumberFormat = 'General'
$ws.Range('E24').Style = 'Normal'

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.48'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.42%  '
$ws.Range('E25').NumberFormat = 'General'
$ws.Range('E25').Style = 'Normal'

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.756'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('D26').Style = 'Normal'

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +16.25%  '
$ws.Range('E26').NumberFormat = 'General'
$ws.Range('E26').Style = 'Normal'

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1235'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +4.51%  '
$ws.Range('E27').NumberFormat = 'General'
$ws.Range('E27').Style = 'Normal'

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.391'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('D28').Style = 'Normal'

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.89%  '
$ws.Range('E28').NumberFormat = 'General'
$ws.Range('E28').Style = 'Normal'

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +4.31%  '
$ws.Range('E29').NumberFormat = 'General'
$ws.Range('E29').Style = 'Normal'

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05489'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = 'Normal'

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.47%  '
$ws.Range('E30').NumberFormat = 'General'
$ws.Range('E30').Style = 'Normal'

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.76%  '
$ws.Range('E31').NumberFormat = 'General'
$ws.Range('E31').Style = 'Normal'

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.26%  '
$ws.Range('E32').NumberFormat = 'General'
$ws.Range('E32').Style = 'Normal'

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +3.40%  '
$ws.Range('E33').NumberFormat = 'General'
$ws.Range('E33').Style = 'Normal'

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.863'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('D35').Style = 'Normal'

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +3.01%  '
$ws.Range('E35').NumberFormat = 'General'
$ws.Range('E35').Style = 'Normal'

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9643'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('D36').Style = 'Normal'

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.43%  '
$ws.Range('E36').NumberFormat = 'General'
$ws.Range('E36').Style = 'Normal'

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.427'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D37').Style = 'Normal'

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.42%  '
$ws.Range('E37').NumberFormat = 'General'
$ws.Range('E37').Style = 'Normal'

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5968'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('D38').Style = 'Normal'

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +6.64%  '
$ws.Range('E38').NumberFormat = 'General'
$ws.Range('E38').Style = 'Normal'

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01649'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('D39').Style = 'Normal'

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +4.95%  '
$ws.Range('E39').NumberFormat = 'General'
$ws.Range('E39').Style = 'Normal'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.902'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('D40').Style = 'Normal'

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.10%  '
$ws.Range('E40').NumberFormat = 'General'
$ws.Range('E40').Style = 'Normal'

$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('B41').NumberFormat = 'General'
$ws.Range('B41').Style = 'Normal'

$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('C41').NumberFormat = 'General'
$ws.Range('C41').Style = 'Normal'

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8523'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('D41').Style = 'Normal'

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.24%  '
$ws.Range('E41').NumberFormat = 'General'
$ws.Range('E41').Style = 'Normal'

$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'Maker'
$ws.Range('B42').NumberFormat = 'General'
$ws.Range('B42').Style = 'Normal'

$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('C42').NumberFormat = 'General'
$ws.Range('C42').Style = 'Normal'

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.054.03'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = 'Normal'

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.69%  '
$ws.Range('E42').NumberFormat = 'General'
$ws.Range('E42').Style = 'Normal'

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.14%  '
$ws.Range('E43').NumberFormat = 'General'
$ws.Range('E43').Style = 'Normal'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.37'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('D44').Style = 'Normal'

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.870.36'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('D45').Style = 'Normal'

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +4.73%  '
$ws.Range('E45').NumberFormat = 'General'
$ws.Range('E45').Style = 'Normal'

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +4.03%  '
$ws.Range('E46').NumberFormat = 'General'
$ws.Range('E46').Style = 'Normal'

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '58.84'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('D47').Style = 'Normal'

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.68%  '
$ws.Range('E47').NumberFormat = 'General'
$ws.Range('E47').Style = 'Normal'

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4476'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.71%  '
$ws.Range('E48').NumberFormat = 'General'
$ws.Range('E48').Style = 'Normal'

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.208'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('D49').Style = 'Normal'

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.87%  '
$ws.Range('E49').NumberFormat = 'General'
$ws.Range('E49').Style = 'Normal'

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.003'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = 'Normal'

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05241'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.54%  '
$ws.Range('E51').NumberFormat = 'General'
$ws.Range('E51').Style = 'Normal'
